# WSJT-X Quick Start Guide - text corrections
# Applies the changes described in the commit: fixes missing spaces,
# hyphenates "sub-mode(s)", reworks a couple of sentences, and corrects
# a scaling error (207 -> 206) in the JT9 spec paragraph.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText, $matchCase = $true) {
    $d.Content.Find.Execute($findText, $matchCase, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Missing space: "modesJT65" -> "modes JT65"
Replace-Text "popular modesJT65 and JT4" "popular modes JT65 and JT4"

# 2. "Submodes with" -> "Sub-modes with"
Replace-Text "Submodes with" "Sub-modes with"

# 3. "The slowest submode, JT9-30" -> "The slowest sub-mode, JT9-30"
Replace-Text "The slowest submode, JT9-30" "The slowest sub-mode, JT9-30"

# 4. "JT9 submode, the" -> "JT9 sub-mode, the"
Replace-Text "JT9 submode, the" "JT9 sub-mode, the"

# 5. "... for submodes JT9-1 and JT9-2 ..." -> "... for sub-modes JT9-1 and JT9-2 ..."
Replace-Text "1000–2000 Hz for submodes JT9-1 and JT9-2" "1000–2000 Hz for sub-modes JT9-1 and JT9-2"

# 6. "Double-click to set QSO Frequency" -> "Double-clicking will set QSO Frequency"
Replace-Text "Double-click to set QSO Frequency" "Double-clicking will set QSO Frequency"

# 7. "... no particular effort has been put into the decoder's handling ..."
#    -> "... no particular effort has been put into optimizing the decoder's handling ..."
Replace-Text "put into the decoder" "put into optimizing the decoder"

# 8. Technical paragraph fixes (scaling error correction + missing-space fixes)
Replace-Text "Error control coding(ECC)" "Error control coding (ECC)"
Replace-Text "K=32,rate r=1/2" "K=32, rate r=1/2"
Replace-Text "9-FSK: 8tones for data" "9-FSK: 8 tones are used for data"
Replace-Text "Sixteen symbol intervals areused" "Sixteen symbol intervals are used"
Replace-Text "a transmission requires a total of 207/3+ 16 = 85 channel symbols" "a transmission requires a total of  206 / 3 + 16 = 85 (rounded up) channel symbols"
Replace-Text "approximately(TRperiod-8)/85" "approximately (TRperiod - 8) / 85"
Replace-Text "in seconds.Exact symbol" "in seconds.  Exact symbol"
Replace-Text "nsps, the number of samplesper symbol" "nsps, the number of samples per symbol"
Replace-Text "no primefactor greater than 7" "no prime factor greater than 7"
Replace-Text "FFTs.  Tonespacing" "FFTs.  Tone spacing"
Replace-Text "df=1/tsym=12000/nsps, equal to thekeying rate" "df = 1 / tsym = 12000 / nsps, equal to the keying rate"
Replace-Text "is 9*df.  The generatedsignal has continuous phase" "is 9 × df.  The generated signal has continuous phase"

# 9. "the followingtable" -> "the following table"
Replace-Text "the followingtable" "the following table"

# 10. Move the "_GoBack" bookmark from its old location (after "to a suitable")
#     to the new location (after "for sub-" before "modes JT9-1").
$d.Bookmarks("_GoBack").Delete()
$r = $d.Content
$ok = $r.Find.Execute("for sub-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $r.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r) | Out-Null
}

Write-Output "done"
